# debugging mac_lane int mode
# - Adds a "row total" column BM for rows 3 and 4 (shared formula BM3:BM4)
# - Adds two new data rows (6 and 7), each a duplicate half of row 3's data
#   (row 6 = AG:BL slice, row 7 = A:AF slice), plus BM totals for rows 6-10
#   (shared formula BM6:BM10, with rows 8-10 empty so their totals are 0)
# - Updates the selected cell shown in the worksheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row-sum column for rows 3 and 4 (shared formula across BM3:BM4) ---
$ws.Range("BM3:BM4").Formula = "=SUM(A3:BL3)"

# --- New row 6: second half (AG:BL) of row 3's values ---
$row6Vals = @(145,140,123,-21,-247,228,9,-221,192,123,-55,213,-178,-84,-171,-131,-126,-176,-198,-145,91,-24,-166,-126,35,-47,203,147,-32,58,-217,208)
$row6Arr = New-Object 'object[,]' 1,$row6Vals.Length
for ($i = 0; $i -lt $row6Vals.Length; $i++) { $row6Arr[0,$i] = $row6Vals[$i] }
$ws.Range("AG6:BL6").Value = $row6Arr

# --- New row 7: first half (A:AF) of row 3's values ---
$row7Vals = @(134,-216,-86,140,196,-73,-141,93,-237,60,183,-229,192,-216,-214,-14,-23,225,7,148,99,-163,133,138,-135,218,96,137,246,-6,-20,-214)
$row7Arr = New-Object 'object[,]' 1,$row7Vals.Length
for ($i = 0; $i -lt $row7Vals.Length; $i++) { $row7Arr[0,$i] = $row7Vals[$i] }
$ws.Range("A7:AF7").Value = $row7Arr

# --- New row-sum column for rows 6 through 10 (shared formula across BM6:BM10) ---
$ws.Range("BM6:BM10").Formula = "=SUM(A6:BL6)"

# --- Update the view's active selection ---
$ws.Range("BK19").Select()
